$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for rows 2-108
$ws.Range("F2").Value = "2021-10-05 13:41:49.092127"
$ws.Range("F3").Value = "2021-10-05 13:41:49.092140"
$ws.Range("F4").Value = "2021-10-05 13:41:49.092144"
$ws.Range("F5").Value = "2021-10-05 13:41:49.092147"
$ws.Range("F6").Value = "2021-10-05 13:41:49.092151"
$ws.Range("F7").Value = "2021-10-05 13:41:49.092154"
$ws.Range("F8").Value = "2021-10-05 13:41:49.092157"
$ws.Range("F9").Value = "2021-10-05 13:41:49.092160"
$ws.Range("F10").Value = "2021-10-05 13:41:49.092163"
$ws.Range("F11").Value = "2021-10-05 13:41:49.092166"
$ws.Range("F12").Value = "2021-10-05 13:41:49.092169"
$ws.Range("F13").Value = "2021-10-05 13:41:49.092172"
$ws.Range("F14").Value = "2021-10-05 13:41:49.092175"
$ws.Range("F15").Value = "2021-10-05 13:41:49.092178"
$ws.Range("F16").Value = "2021-10-05 13:41:49.092181"
$ws.Range("F17").Value = "2021-10-05 13:41:49.092184"
$ws.Range("F18").Value = "2021-10-05 13:41:49.092187"
$ws.Range("F19").Value = "2021-10-05 13:41:49.092191"
$ws.Range("F20").Value = "2021-10-05 13:41:49.092194"
$ws.Range("F21").Value = "2021-10-05 13:41:49.092196"
$ws.Range("F22").Value = "2021-10-05 13:41:49.092199"
$ws.Range("F23").Value = "2021-10-05 13:41:49.092202"
$ws.Range("F24").Value = "2021-10-05 13:41:49.092205"
$ws.Range("F25").Value = "2021-10-05 13:41:49.092208"
$ws.Range("F26").Value = "2021-10-05 13:41:49.092212"
$ws.Range("F27").Value = "2021-10-05 13:41:49.092215"
$ws.Range("F28").Value = "2021-10-05 13:41:49.092218"
$ws.Range("F29").Value = "2021-10-05 13:41:49.092221"
$ws.Range("F30").Value = "2021-10-05 13:41:49.092224"
$ws.Range("F31").Value = "2021-10-05 13:41:49.092227"
$ws.Range("F32").Value = "2021-10-05 13:41:49.092230"
$ws.Range("F33").Value = "2021-10-05 13:41:49.092233"
$ws.Range("F34").Value = "2021-10-05 13:41:49.092237"
$ws.Range("F35").Value = "2021-10-05 13:41:49.092240"
$ws.Range("F36").Value = "2021-10-05 13:41:49.092243"
$ws.Range("F37").Value = "2021-10-05 13:41:49.092246"
$ws.Range("F38").Value = "2021-10-05 13:41:49.092249"
$ws.Range("F39").Value = "2021-10-05 13:41:49.092252"
$ws.Range("F40").Value = "2021-10-05 13:41:49.092255"
$ws.Range("F41").Value = "2021-10-05 13:41:49.092258"
$ws.Range("F42").Value = "2021-10-05 13:41:49.092261"
$ws.Range("F43").Value = "2021-10-05 13:41:49.092264"
$ws.Range("F44").Value = "2021-10-05 13:41:49.092267"
$ws.Range("F45").Value = "2021-10-05 13:41:49.092270"
$ws.Range("F46").Value = "2021-10-05 13:41:49.092273"
$ws.Range("F47").Value = "2021-10-05 13:41:49.092276"
$ws.Range("F48").Value = "2021-10-05 13:41:49.092279"
$ws.Range("F49").Value = "2021-10-05 13:41:49.092282"
$ws.Range("F50").Value = "2021-10-05 13:41:49.092285"
$ws.Range("F51").Value = "2021-10-05 13:41:49.092288"
$ws.Range("F52").Value = "2021-10-05 13:41:49.092291"
$ws.Range("F53").Value = "2021-10-05 13:41:49.092294"
$ws.Range("F54").Value = "2021-10-05 13:41:49.092297"
$ws.Range("F55").Value = "2021-10-05 13:41:49.092301"
$ws.Range("F56").Value = "2021-10-05 13:41:49.092304"
$ws.Range("F57").Value = "2021-10-05 13:41:49.092307"
$ws.Range("F58").Value = "2021-10-05 13:41:49.092310"
$ws.Range("F59").Value = "2021-10-05 13:41:49.092313"
$ws.Range("F60").Value = "2021-10-05 13:41:49.092316"
$ws.Range("F61").Value = "2021-10-05 13:41:49.092319"
$ws.Range("F62").Value = "2021-10-05 13:41:49.092322"
$ws.Range("F63").Value = "2021-10-05 13:41:49.092325"
$ws.Range("F64").Value = "2021-10-05 13:41:49.092328"
$ws.Range("F65").Value = "2021-10-05 13:41:49.092331"
$ws.Range("F66").Value = "2021-10-05 13:41:49.092335"
$ws.Range("F67").Value = "2021-10-05 13:41:49.092339"
$ws.Range("F68").Value = "2021-10-05 13:41:49.092342"
$ws.Range("F69").Value = "2021-10-05 13:41:49.092345"
$ws.Range("F70").Value = "2021-10-05 13:41:49.092348"
$ws.Range("F71").Value = "2021-10-05 13:41:49.092351"
$ws.Range("F72").Value = "2021-10-05 13:41:49.092353"
$ws.Range("F73").Value = "2021-10-05 13:41:49.092357"
$ws.Range("F74").Value = "2021-10-05 13:41:49.092360"
$ws.Range("F75").Value = "2021-10-05 13:41:49.092363"
$ws.Range("F76").Value = "2021-10-05 13:41:49.092366"
$ws.Range("F77").Value = "2021-10-05 13:41:49.092369"
$ws.Range("F78").Value = "2021-10-05 13:41:49.092374"
$ws.Range("F79").Value = "2021-10-05 13:41:49.092377"
$ws.Range("F80").Value = "2021-10-05 13:41:49.092380"
$ws.Range("F81").Value = "2021-10-05 13:41:49.092383"
$ws.Range("F82").Value = "2021-10-05 13:41:49.092386"
$ws.Range("F83").Value = "2021-10-05 13:41:49.092389"
$ws.Range("F84").Value = "2021-10-05 13:41:49.092392"
$ws.Range("F85").Value = "2021-10-05 13:41:49.092395"
$ws.Range("F86").Value = "2021-10-05 13:41:49.092398"
$ws.Range("F87").Value = "2021-10-05 13:41:49.092401"
$ws.Range("F88").Value = "2021-10-05 13:41:49.092404"
$ws.Range("F89").Value = "2021-10-05 13:41:49.092407"
$ws.Range("F90").Value = "2021-10-05 13:41:49.092410"
$ws.Range("F91").Value = "2021-10-05 13:41:49.092413"
$ws.Range("F92").Value = "2021-10-05 13:41:49.092416"
$ws.Range("F93").Value = "2021-10-05 13:41:49.092419"
$ws.Range("F94").Value = "2021-10-05 13:41:49.092423"
$ws.Range("F95").Value = "2021-10-05 13:41:49.092427"
$ws.Range("F96").Value = "2021-10-05 13:41:49.092430"
$ws.Range("F97").Value = "2021-10-05 13:41:49.092433"
$ws.Range("F98").Value = "2021-10-05 13:41:49.092436"
$ws.Range("F99").Value = "2021-10-05 13:41:49.092439"
$ws.Range("F100").Value = "2021-10-05 13:41:49.092442"
$ws.Range("F101").Value = "2021-10-05 13:41:49.092445"
$ws.Range("F102").Value = "2021-10-05 13:41:49.092448"
$ws.Range("F103").Value = "2021-10-05 13:41:49.092451"
$ws.Range("F104").Value = "2021-10-05 13:41:49.092454"
$ws.Range("F105").Value = "2021-10-05 13:41:49.092457"
$ws.Range("F106").Value = "2021-10-05 13:41:49.092460"
$ws.Range("F107").Value = "2021-10-05 13:41:49.092463"
$ws.Range("F108").Value = "2021-10-05 13:41:49.092466"
